$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2796.6667
$ws.Range("J2").Value = 4673.875
$ws.Range("L2").Value = 4673.875
$ws.Range("N2").Value = -4899.875

$ws.Range("H6").Value = 1650.0769
$ws.Range("I6").Value = 2301.3333
$ws.Range("J6").Value = 184.75
$ws.Range("K6").Value = 6903.999899999999
$ws.Range("L6").Value = 554.25
$ws.Range("M6").Value = -6791.999899999999
$ws.Range("N6").Value = -778.25

$ws.Range("H43").Value = 877.4
$ws.Range("I43").Value = 797
$ws.Range("J43").Value = 998
$ws.Range("K43").Value = 797
$ws.Range("L43").Value = 998
$ws.Range("M43").Value = -728
$ws.Range("N43").Value = -1136

$ws.Range("H55").Value = 550
$ws.Range("J55").Value = 650
$ws.Range("L55").Value = 650
$ws.Range("N55").Value = -1078

$ws.Range("H58").Value = 928.8125
$ws.Range("J58").Value = 1666.6666
$ws.Range("L58").Value = 4999.9998
$ws.Range("N58").Value = -5299.9998

$ws.Range("I70").Value = 1325.2
$ws.Range("J70").Value = 3920.7778
$ws.Range("K70").Value = 3975.6
$ws.Range("L70").Value = 11762.3334
$ws.Range("M70").Value = -3705.6
$ws.Range("N70").Value = -12302.3334

$ws.Range("I73").Value = 1325.2
$ws.Range("J73").Value = 3920.7778
$ws.Range("K73").Value = 3975.6
$ws.Range("L73").Value = 11762.3334
$ws.Range("M73").Value = -3039.6
$ws.Range("N73").Value = -13634.3334

$ws.Range("H107").Value = 2783.5557
$ws.Range("J107").Value = 10002.5
$ws.Range("L107").Value = 10002.5
$ws.Range("N107").Value = -13842.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 1097.5714
$ws.Range("I35").Value = 589.5
$ws.Range("K35").Value = 589.5
$ws.Range("M35").Value = -183.5

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 25743.857

$ws.Range("H85").Value = 25743.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1928.5385
$ws.Range("I58").Value = 1928.5385
$ws.Range("K58").Value = 1928.5385
$ws.Range("M58").Value = -1725.5385

$ws.Range("H134").Value = 3301.4736
$ws.Range("I134").Value = 3159.4119
$ws.Range("K134").Value = 9478.235700000001
$ws.Range("M134").Value = -6943.235700000001

$ws.Range("H136").Value = 1928.5385
$ws.Range("I136").Value = 1928.5385
$ws.Range("K136").Value = 5785.6155
$ws.Range("M136").Value = -3235.6155

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws.Range("H24").Value = 720
$ws.Range("I24").Value = 1035
$ws.Range("J24").Value = 90
$ws.Range("K24").Value = 3105
$ws.Range("L24").Value = 270
$ws.Range("M24").Value = -2875
$ws.Range("N24").Value = -730

$ws.Range("H34").Value = 506.83334
$ws.Range("I34").Value = 414
$ws.Range("J34").Value = 599.6667
$ws.Range("K34").Value = 1242
$ws.Range("L34").Value = 1799.0001
$ws.Range("M34").Value = -1158
$ws.Range("N34").Value = -1967.0001

$ws.Range("H36").Value = 287.66666
$ws.Range("I36").Value = 305.75
$ws.Range("K36").Value = 917.25
$ws.Range("M36").Value = -748.25

$ws.Range("H41").Value = 1334.3334
$ws.Range("J41").Value = 1334.3334
$ws.Range("L41").Value = 4003.0002
$ws.Range("N41").Value = -4679.0002

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 788
$ws.Range("I134").Value = 788
$ws.Range("K134").Value = 2364
$ws.Range("M134").Value = 2706

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5080.6665
$ws.Range("I80").Value = 2995
$ws.Range("J80").Value = 5497.8
$ws.Range("K80").Value = 2995
$ws.Range("L80").Value = 5497.8
$ws.Range("M80").Value = -1997
$ws.Range("N80").Value = -7493.8

$ws.Range("H83").Value = 5080.6665
$ws.Range("I83").Value = 2995
$ws.Range("J83").Value = 5497.8
$ws.Range("K83").Value = 14975
$ws.Range("L83").Value = 27489
$ws.Range("M83").Value = -9983
$ws.Range("N83").Value = -37473

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1816
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 1816
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H132").Value = 6556.2856
$ws.Range("I132").Value = 2197
$ws.Range("K132").Value = 6591
$ws.Range("M132").Value = -4061

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 8500
$ws.Range("I18").Value = 5000
$ws.Range("J18").Value = 12000
$ws.Range("K18").Value = 5000
$ws.Range("L18").Value = 12000
$ws.Range("M18").Value = -4827
$ws.Range("N18").Value = -12346

$ws.Range("H54").Value = 20750
$ws.Range("I54").Value = 11500
$ws.Range("K54").Value = 11500
$ws.Range("M54").Value = -10980

$ws.Range("H62").Value = 4724.875
$ws.Range("I62").Value = 4959.4
$ws.Range("J62").Value = 4334
$ws.Range("K62").Value = 4959.4
$ws.Range("L62").Value = 4334
$ws.Range("M62").Value = -4335.4
$ws.Range("N62").Value = -5582

$ws.Range("H65").Value = 4724.875
$ws.Range("I65").Value = 4959.4
$ws.Range("J65").Value = 4334
$ws.Range("K65").Value = 24797
$ws.Range("L65").Value = 21670
$ws.Range("M65").Value = -21677
$ws.Range("N65").Value = -27910

$ws.Range("H81").Value = 1966.6666
$ws.Range("I81").Value = 450
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 900
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = 161
$ws.Range("N81").Value = -12122

$ws.Range("H84").Value = 1966.6666
$ws.Range("I84").Value = 450
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 4500
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = 804
$ws.Range("N84").Value = -60608

$ws.Range("H136").Value = 1760.3214
$ws.Range("I136").Value = 1367.56
$ws.Range("J136").Value = 5033.3335
$ws.Range("K136").Value = 4102.68
$ws.Range("L136").Value = 15100.0005
$ws.Range("M136").Value = -1552.68
$ws.Range("N136").Value = -20200.0005
